$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Convert the string "true" in H2/H3 to the real boolean TRUE, and apply the
# same formatting (wrap text, general number format) down through H8 on the
# "survey" sheet (this clears the old numFmtId=49 text format).
$rng = $survey.Range("H2:H8")
$rng.NumberFormat = "General"
$rng.WrapText = $true

$survey.Range("H2").Value = $true
$survey.Range("H3").Value = $true

# Selection / active sheet: survey tab becomes selected (was settings).
$survey.Range("H4").Select()
$survey.Activate()

$wb.Worksheets.Item(1).Activate()
